# ---------------------------------------------------------------------------
# Applies the "Add files via upload" commit to 2_4_DEQ.xlsx:
#   * Tweaks two existing questions on sheet "16_" (wording + row height).
#   * Appends two brand-new question sheets ("17_" and "18_") at the end of
#     the workbook, each holding a small multiple-choice-style question grid
#     (prompt / wrong answer / right answer + feedback), mirroring the look
#     of the other question sheets already in the workbook.
#   * Leaves the newly added "18_" sheet active/selected, like the original
#     author ended up doing after appending the new content.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update sheet "16_" (the existing last question sheet)
# ---------------------------------------------------------------------------
$ws16 = $wb.Worksheets.Item("16_")

# The prompt grew a little longer ("...looks like the plot in the change...")
# so the header row needed to be a bit taller.
$ws16.Rows(1).RowHeight = 120

$ws16.Range("A1").Value = "The plot of the voltage at the capacitor plate looks like the plot in the change in the temperature of the milk returning to the environmental temperature.  How might the two things be similar?  (Mark all that are true)"
$ws16.Range("A4").Value = "Milk is actually a special kind of capacitor"

# The selection on this sheet also moved down one row (E3 -> E4) once the
# author finished editing it.
$ws16.Range("E4").Select()

# ---------------------------------------------------------------------------
# 2) Add new sheet "17_" after the last sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws17 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws17.Name = "17_"

$ws17.Columns("A").ColumnWidth = 25.307291666666668
$ws17.Columns("C").ColumnWidth = 47.022135416666664

$ws17.Range("A1:H19").WrapText = $true

$ws17.Range("A1").Value = "An RC circuit ""filters"" a high frequency signal because it can't move electrons fast enough to keep up with the input signal.   What would the effect of increasing the resistance of the resistor be in such a filter?"

$ws17.Range("A2").Value = "It would allow electrons to move more quickly and so make the V_out more responsive to the input signal"
$ws17.Range("B2").Value = "N"

$ws17.Range("A3").Value = "It would slow down the movement of electrons and so make V_out less responsive to the input signal"
$ws17.Range("B3").Value = "Y"
$ws17.Range("C3").Value = "Yep!  A higher resistance would lower the current (given the same voltage difference), and it would take longer for V_out to approach the voltage at V_in"

$ws17.Rows(1).RowHeight = 120
$ws17.Rows(2).RowHeight = 75
$ws17.Rows(3).RowHeight = 75

$ws17.Range("A1:C3").Select()

# ---------------------------------------------------------------------------
# 3) Add new sheet "18_" after "17_"
# ---------------------------------------------------------------------------
$ws18 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws17)
$ws18.Name = "18_"

$ws18.Columns("A").ColumnWidth = 37.307291666666664
$ws18.Columns("C").ColumnWidth = 40.451822916666664

$ws18.Range("A1:C3").WrapText = $true

$ws18.Range("A1").Value = "An RC circuit ""filters"" a high frequency signal because it can't move electrons fast enough to keep up with the input signal.   A larger resistor would make the circuit react even more slowly.   What would this do the ""cutoff"" frequency of the filter?"

$ws18.Range("A2").Value = "The cutoff frequency would go up: only really high frequencies would be filtered out"
$ws18.Range("B2").Value = "N"

$ws18.Range("A3").Value = "The cutoff frequency would go down: the filter would remove more mid-range frequencies with the larger resistor"
$ws18.Range("B3").Value = "Y"
$ws18.Range("C3").Value = "Yep!  The higher resistor slows the response time of V_out to V_in, so you'd need an even lower frequency input signal in order to get through the filter."

$ws18.Rows(1).RowHeight = 105
$ws18.Rows(2).RowHeight = 45
$ws18.Rows(3).RowHeight = 60

# Final view state: "18_" tab active/selected, cursor parked just below the
# question grid (C4) -- matches where the author left off after typing.
$ws18.Activate()
$ws18.Range("C4").Select()
